# Apply odds/market updates scraped for the 2025-02-08 FlashScore weekly workbook.
# Each row below corresponds to one match; only the cells that moved are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("Q4").Value = 1.88
$ws.Range("R4").Value = 2.02

# Row 22
$ws.Range("G22").Value = 2.3
$ws.Range("H22").Value = 2.63
$ws.Range("J22").Value = 3.25
$ws.Range("K22").Value = 1.73
$ws.Range("M22").Value = 1.18
$ws.Range("N22").Value = 4.5
$ws.Range("O22").Value = 1.73
$ws.Range("P22").Value = 2
$ws.Range("Q22").Value = 3.6
$ws.Range("R22").Value = 1.29
$ws.Range("S22").Value = 8
$ws.Range("T22").Value = 1.08
$ws.Range("U22").Value = 1.78
$ws.Range("V22").Value = 2.03
$ws.Range("W22").Value = 2.63
$ws.Range("X22").Value = 1.44
$ws.Range("AB22").Value = 23
$ws.Range("AC22").Value = 29
$ws.Range("AD22").Value = 51
$ws.Range("AE22").Value = 4.5
$ws.Range("AH22").Value = 126
$ws.Range("AJ22").Value = 7
$ws.Range("AL22").Value = 17
$ws.Range("AO22").Value = 67
$ws.Range("AP22").Value = 6.2
$ws.Range("AQ22").Value = 1.12

# Row 23
$ws.Range("G23").Value = 1.73
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 5.25
$ws.Range("J23").Value = 2.5
$ws.Range("L23").Value = 6
$ws.Range("M23").Value = 1.11
$ws.Range("N23").Value = 6.5
$ws.Range("W23").Value = 2.38
$ws.Range("X23").Value = 1.53
$ws.Range("Z23").Value = 6.5
$ws.Range("AB23").Value = 13
$ws.Range("AC23").Value = 17
$ws.Range("AF23").Value = 7
$ws.Range("AG23").Value = 23
$ws.Range("AJ23").Value = 10
$ws.Range("AL23").Value = 19
$ws.Range("AN23").Value = 51

# Row 24
$ws.Range("G24").Value = 1.73
$ws.Range("H24").Value = 3.7
$ws.Range("J24").Value = 2.4
$ws.Range("K24").Value = 2.05
$ws.Range("Q24").Value = 2.25
$ws.Range("R24").Value = 1.62
$ws.Range("AK24").Value = 23
$ws.Range("AP24").Value = 3.6
$ws.Range("AQ24").Value = 1.3

# Row 25
$ws.Range("G25").Value = 2.3
$ws.Range("H25").Value = 3.2
$ws.Range("I25").Value = 3.1
$ws.Range("J25").Value = 3.1
$ws.Range("L25").Value = 3.75
$ws.Range("N25").Value = 8.5
$ws.Range("W25").Value = 1.91
$ws.Range("X25").Value = 1.91
$ws.Range("Z25").Value = 11
$ws.Range("AA25").Value = 9.5
$ws.Range("AB25").Value = 21
$ws.Range("AC25").Value = 21
$ws.Range("AE25").Value = 8.5
$ws.Range("AF25").Value = 6
$ws.Range("AG25").Value = 15
$ws.Range("AI25").Value = 301
$ws.Range("AJ25").Value = 9
$ws.Range("AN25").Value = 26

# Row 46
$ws.Range("G46").Value = 1.4
$ws.Range("H46").Value = 4.33
$ws.Range("I46").Value = 9.5
$ws.Range("M46").Value = 1.08
$ws.Range("N46").Value = 8
$ws.Range("AE46").Value = 8
$ws.Range("AF46").Value = 9
$ws.Range("AG46").Value = 29
$ws.Range("AH46").Value = 101
$ws.Range("AJ46").Value = 17
$ws.Range("AL46").Value = 29
$ws.Range("AM46").Value = 126
$ws.Range("AN46").Value = 81
$ws.Range("AO46").Value = 81

# Row 47
$ws.Range("G47").Value = 2.88
$ws.Range("H47").Value = 2.55
$ws.Range("I47").Value = 3
$ws.Range("J47").Value = 3.6
$ws.Range("K47").Value = 1.91
$ws.Range("L47").Value = 3.75
$ws.Range("M47").Value = 1.13
$ws.Range("N47").Value = 6
$ws.Range("O47").Value = 1.53
$ws.Range("P47").Value = 2.38
$ws.Range("Q47").Value = 2.7
$ws.Range("R47").Value = 1.44
$ws.Range("S47").Value = 5.5
$ws.Range("T47").Value = 1.14
$ws.Range("U47").Value = 1.62
$ws.Range("V47").Value = 2.2
$ws.Range("W47").Value = 2.1
$ws.Range("X47").Value = 1.67
$ws.Range("Y47").Value = 7
$ws.Range("Z47").Value = 12
$ws.Range("AA47").Value = 12
$ws.Range("AB47").Value = 29
$ws.Range("AE47").Value = 5.5
$ws.Range("AG47").Value = 17
$ws.Range("AH47").Value = 67
$ws.Range("AK47").Value = 13
$ws.Range("AL47").Value = 12
$ws.Range("AM47").Value = 29
$ws.Range("AN47").Value = 29
$ws.Range("AO47").Value = 41
$ws.Range("AP47").Value = 4.4
$ws.Range("AQ47").Value = 1.2
$ws.Range("AR47").Value = 2.05
$ws.Range("AS47").Value = 1.8

# Row 48
$ws.Range("G48").Value = 1.6
$ws.Range("H48").Value = 3.7
$ws.Range("I48").Value = 5.75
$ws.Range("J48").Value = 2.3

# Row 108
$ws.Range("G108").Value = 3.2
$ws.Range("H108").Value = 3.6
$ws.Range("J108").Value = 3.5
$ws.Range("Q108").Value = 1.67
$ws.Range("R108").Value = 2.15
$ws.Range("AC108").Value = 21
$ws.Range("AF108").Value = 7.5
$ws.Range("AN108").Value = 17
$ws.Range("AO108").Value = 23

# Row 109
$ws.Range("G109").Value = 1.48
$ws.Range("H109").Value = 4.33
$ws.Range("I109").Value = 6.5
$ws.Range("J109").Value = 2.05
$ws.Range("L109").Value = 6.5
$ws.Range("W109").Value = 2
$ws.Range("X109").Value = 1.75
$ws.Range("AB109").Value = 10
$ws.Range("AF109").Value = 8
$ws.Range("AI109").Value = 401

# Row 111
$ws.Range("G111").Value = 1.62
$ws.Range("H111").Value = 3.6
$ws.Range("I111").Value = 5.1
$ws.Range("J111").Value = 2.22
$ws.Range("K111").Value = 2.12
$ws.Range("L111").Value = 5.1
$ws.Range("O111").Value = 1.26
$ws.Range("P111").Value = 3.15
$ws.Range("Q111").Value = 1.78
$ws.Range("R111").Value = 1.82
$ws.Range("S111").Value = 2.82
$ws.Range("T111").Value = 1.32
$ws.Range("U111").Value = 1.39
$ws.Range("V111").Value = 2.57
$ws.Range("W111").Value = 1.78
$ws.Range("X111").Value = 1.83
$ws.Range("Z111").Value = 7.6
$ws.Range("AB111").Value = 12.5
$ws.Range("AC111").Value = 13
$ws.Range("AE111").Value = 10.25
$ws.Range("AF111").Value = 7.1
$ws.Range("AG111").Value = 15.5
$ws.Range("AH111").Value = 70
$ws.Range("AI111").Value = 600
$ws.Range("AJ111").Value = 14
$ws.Range("AK111").Value = 30
$ws.Range("AL111").Value = 16
$ws.Range("AM111").Value = 90
$ws.Range("AN111").Value = 50
$ws.Range("AO111").Value = 50

# Row 121
$ws.Range("M121").Value = 1.05
$ws.Range("N121").Value = 11
$ws.Range("Q121").Value = 1.88
$ws.Range("R121").Value = 1.98
$ws.Range("S121").Value = 3.25
$ws.Range("T121").Value = 1.33
$ws.Range("AP121").Value = 2.43
$ws.Range("AQ121").Value = 1.57
$ws.Range("AR121").Value = 1.41
$ws.Range("AS121").Value = 2.95

# Row 170
$ws.Range("G170").Value = 1.72
$ws.Range("H170").Value = 2.92
$ws.Range("I170").Value = 5.9
$ws.Range("J170").Value = 2.3
$ws.Range("K170").Value = 1.93
$ws.Range("L170").Value = 5.9
$ws.Range("O170").Value = 1.44
$ws.Range("P170").Value = 2.4
$ws.Range("Q170").Value = 2.27
$ws.Range("R170").Value = 1.5
$ws.Range("S170").Value = 3.85
$ws.Range("W170").Value = 2.05
$ws.Range("X170").Value = 1.6
$ws.Range("Y170").Value = 5
$ws.Range("Z170").Value = 7
$ws.Range("AA170").Value = 8.25
$ws.Range("AB170").Value = 14
$ws.Range("AC170").Value = 16.5
$ws.Range("AE170").Value = 6.3
$ws.Range("AF170").Value = 6
$ws.Range("AG170").Value = 18
$ws.Range("AH170").Value = 110
$ws.Range("AJ170").Value = 12
$ws.Range("AK170").Value = 35
$ws.Range("AM170").Value = 150
$ws.Range("AN170").Value = 80
$ws.Range("AO170").Value = 80

# Row 171
$ws.Range("G171").Value = 2.65
$ws.Range("H171").Value = 2.8
$ws.Range("I171").Value = 2.82
$ws.Range("J171").Value = 3.4
$ws.Range("K171").Value = 1.85
$ws.Range("L171").Value = 3.5
$ws.Range("M171").Value = 1.11
$ws.Range("N171").Value = 6.2
$ws.Range("O171").Value = 1.52
$ws.Range("P171").Value = 2.2
$ws.Range("Q171").Value = 2.47
$ws.Range("R171").Value = 1.42
$ws.Range("S171").Value = 4.3
$ws.Range("T171").Value = 1.14
$ws.Range("U171").Value = 1.57
$ws.Range("V171").Value = 2.12
$ws.Range("W171").Value = 2.05
$ws.Range("X171").Value = 1.6
$ws.Range("Y171").Value = 6.1
$ws.Range("AA171").Value = 10.75
$ws.Range("AB171").Value = 32
$ws.Range("AC171").Value = 29
$ws.Range("AE171").Value = 6
$ws.Range("AF171").Value = 5.6
$ws.Range("AK171").Value = 13
